# Update the "Enterprises density (per 1000 people)" row (row 13) values
# for Micro / SMEs / MSMEs from 32.8 / 8.7 / 41.5 to 32.84 / 8.66 / 41.51.
# The source values are stored as text (shared strings), so we force the
# range to Text format before assigning, then restore the cell style back
# to Normal/General so the cells keep their original look & formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

$rng = $ws.Range("B13:D13")
$rng.NumberFormat = "@"

$ws.Range("B13").Value = "32.84"
$ws.Range("C13").Value = "8.66"
$ws.Range("D13").Value = "41.51"

$rng.Style = "Normal"
